$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace "Bermagui NSW" exposure site with "Lakes Entrance" site
$ws.Range("A2").Value = "Lakes Entrance"
$ws.Range("B2").Value = "Chants Summer Carnival - Footbridge, Lakes Entrance VIC 3909"
$ws.Range("C2").Value = "29/12/2020 7:00pm-9:30pm"
$ws.Range("D2").Value = "Case attended carnival"
$ws.Range("E2").Value = "new"

# Row 3: replace "Eden NSW" exposure site with "Melbourne" site (old entry)
$ws.Range("A3").Value = "Melbourne"
$ws.Range("B3").Value = "European Bier Cafe City  120 Exhibition Street Melbourne VIC 3000"
$ws.Range("C3").Value = "21/12/2020 8:00pm-9:30pm"
$ws.Range("D3").Value = "Case attended cafe"
$ws.Range("E3").Value = "old"

# Row 4: new entry, same site as row 3, updated date, marked "new"
$ws.Range("A4").Value = "Melbourne"
$ws.Range("B4").Value = "European Bier Cafe City  120 Exhibition Street Melbourne VIC 3000"
$ws.Range("C4").Value = "28/12/2020 8:00pm-9:30pm"
$ws.Range("D4").Value = "Case attended cafe"
$ws.Range("E4").Value = "new"

# Row 5: new site "Southbank", old entry
$ws.Range("A5").Value = "Southbank"
$ws.Range("B5").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C5").Value = "23/112/2020 1:00pm-1:30pm"
$ws.Range("D5").Value = "Case attended restaurant"
$ws.Range("E5").Value = "old"

# Row 6: same site as row 5, corrected date, marked "new"
$ws.Range("A6").Value = "Southbank"
$ws.Range("B6").Value = "Rockpool Bar and Grill, Crown Casino  8 Whiteman Street, Southbank"
$ws.Range("C6").Value = "23/12/2020 1:00pm-1:30pm"
$ws.Range("D6").Value = "Case attended restaurant"
$ws.Range("E6").Value = "new"

# Leave selection on B2, matching the saved workbook view
$ws.Range("B2").Select() | Out-Null
